$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update bias of post_lasso (n=120, m=80) -- cell C2
$ws.Range("C2").Value = 0.0301

# Related recalculated RJIVE bias/RMSE values
$ws.Range("C4").Value = -0.0119
$ws.Range("E4").Value = -0.0064
$ws.Range("F4").Value = -0.0047

$ws.Range("C5").Value = 0.0871
$ws.Range("D5").Value = 0.1215
$ws.Range("E5").Value = 0.0477
$ws.Range("F5").Value = 0.0559

# Update selection to match author's final cursor position
$ws.Range("C6").Select()
